$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates derived from the authoritative diff. Columns containing
# values that could be auto-parsed by Excel as numbers (column D, "Price")
# are forced to remain plain text, matching the original inlineStr cells,
# by temporarily applying a text number format and resetting the style
# back to Normal afterwards (so no stray style index is left on the cell).
function Set-TextValue($cell, $text) {
    $r = $ws.Range($cell)
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = "Normal"
}

Set-TextValue "D2" "30.520.66"
$ws.Range("E2").Value = "  +0.39%  "
Set-TextValue "D3" "2.135.80"
$ws.Range("E3").Value = "  +1.71%  "
$ws.Range("E4").Value = "  +0.51%  "
Set-TextValue "D5" "352.15"
$ws.Range("E6").Value = "  +0.51%  "
Set-TextValue "D7" "0.5260"
$ws.Range("E7").Value = "  +0.66%  "
Set-TextValue "D8" "0.4565"
$ws.Range("E8").Value = "  -0.07%  "
Set-TextValue "D9" "53.61"
$ws.Range("E9").Value = "  -5.16%  "
Set-TextValue "D10" "0.09156"
$ws.Range("E10").Value = "  +2.59%  "
Set-TextValue "D11" "1.191"
$ws.Range("E11").Value = "  +1.16%  "
Set-TextValue "D12" "25.53"
$ws.Range("E12").Value = "  +5.55%  "
Set-TextValue "D13" "2.138.10"
$ws.Range("E13").Value = "  +1.79%  "
Set-TextValue "D14" "6.893"
$ws.Range("E14").Value = "  +1.01%  "
Set-TextValue "D15" "8.192"
$ws.Range("E15").Value = "  +1.77%  "
Set-TextValue "D16" "101.27"
$ws.Range("E16").Value = "  +4.14%  "
Set-TextValue "D17" "0.00001169"
$ws.Range("E17").Value = "  +1.73%  "
Set-TextValue "D18" "1.009"
$ws.Range("E18").Value = "  +0.54%  "
Set-TextValue "D19" "0.06723"
$ws.Range("E19").Value = "  +1.39%  "
Set-TextValue "D20" "20.44"
$ws.Range("E20").Value = "  +6.63%  "
$ws.Range("E21").Value = "  +0.46%  "
Set-TextValue "D22" "6.380"
$ws.Range("E22").Value = "  +1.19%  "
Set-TextValue "D23" "30.616.56"
$ws.Range("E23").Value = "  +0.46%  "
Set-TextValue "D24" "12.85"
$ws.Range("E24").Value = "  +3.95%  "
Set-TextValue "D25" "2.373"
$ws.Range("E25").Value = "  +0.59%  "
Set-TextValue "D26" "2.386.63"
$ws.Range("E26").Value = "  +1.78%  "
$ws.Range("E27").Value = "  +1.38%  "
Set-TextValue "D28" "2.599"
Set-TextValue "D29" "165.19"
$ws.Range("E29").Value = "  +1.38%  "
Set-TextValue "D30" "135.66"
$ws.Range("E30").Value = "  +1.79%  "
Set-TextValue "D31" "1.226"
$ws.Range("E31").Value = "  +1.32%  "
Set-TextValue "D32" "0.1083"
$ws.Range("E32").Value = "  +1.27%  "
Set-TextValue "D33" "1.714"
$ws.Range("E33").Value = "  +3.42%  "
Set-TextValue "D34" "6.401"
$ws.Range("E34").Value = "  +0.58%  "
Set-TextValue "D35" "4.038"
$ws.Range("E35").Value = "  +2.48%  "
Set-TextValue "D36" "6.153"
$ws.Range("E36").Value = "  +4.06%  "
$ws.Range("E37").Value = "  +1.31%  "
$ws.Range("E38").Value = "  +2.41%  "
Set-TextValue "D39" "0.06982"
$ws.Range("E39").Value = "  +1.82%  "
Set-TextValue "D41" "12.72"
$ws.Range("E41").Value = "  +0.42%  "
Set-TextValue "D42" "0.6985"
$ws.Range("E42").Value = "  +1.66%  "
Set-TextValue "D43" "1.276"
$ws.Range("E43").Value = "  +2.46%  "
Set-TextValue "D44" "14.76"
$ws.Range("E44").Value = "  +5.20%  "
Set-TextValue "D45" "0.6511"
$ws.Range("E45").Value = "  +1.86%  "
$ws.Range("E46").Value = "  +0.95%  "
Set-TextValue "D47" "0.00000000374"
$ws.Range("E47").Value = "  +10.33%  "
Set-TextValue "D48" "3.715"
$ws.Range("E48").Value = "  +1.49%  "
$ws.Range("E49").Value = "  +0.07%  "
Set-TextValue "D50" "83.80"
$ws.Range("E50").Value = "  +0.99%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D51" "0.07287"
$ws.Range("E51").Value = "  +2.22%  "
